$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: RowIndex, D(Fecha), J(Volumen), K(Precio minimo), L(Precio maximo), M(Precio promedio ponderado), N(Unidad de comercializacion), O(Origen), P(Precio $/Kg), Q(Kg o Unidades)
$data = @(
    @(108, 45072, 50, 26000, 26000, 26000, "`$/caja 18 kilos", "Perú", 1444, 18),
    @(109, 45033, 25, 26000, 26000, 26000, "`$/caja 18 kilos", "Perú", 1444, 18),
    @(110, 45062, 5, 26000, 26000, 26000, "`$/caja 18 kilos", "Perú", 1444, 18),
    @(111, 45068, 80, 26000, 26000, 26000, "`$/caja 18 kilos", "Perú", 1444, 18),
    @(112, 44321, 15, 25000, 25000, 25000, "`$/caja 15 kilos granel", "Perú", 1667, 15),
    @(113, 44294, 5, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(114, 44424, 30, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(115, 44316, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(116, 44438, 40, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(117, 44754, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(118, 44790, 15, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(119, 44613, 30, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(120, 44882, 50, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(121, 44819, 100, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(122, 44936, 15, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(123, 44627, 20, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(124, 44329, 40, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(125, 44985, 30, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(126, 44188, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(127, 44637, 30, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(128, 44837, 40, 25000, 25000, 25000, "`$/malla 20 kilos", "Perú", 1250, 20),
    @(129, 44629, 30, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(130, 44810, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(131, 44781, 80, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(132, 45030, 40, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(133, 44617, 20, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(134, 45012, 40, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(135, 44931, 80, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(136, 44665, 40, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(137, 44634, 20, 17000, 17000, 17000, "`$/caja 15 kilos granel", "Perú", 1133, 15),
    @(138, 44634, 30, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(139, 44642, 25, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(140, 44783, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(141, 44893, 40, 22000, 22000, 22000, "`$/malla 20 kilos", "Perú", 1100, 20),
    @(142, 44811, 50, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(143, 44677, 20, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(144, 44881, 40, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(145, 44993, 8, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(146, 44645, 30, 18000, 18000, 18000, "`$/caja 15 kilos granel", "Perú", 1200, 15),
    @(147, 44932, 20, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(148, 44579, 50, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(149, 44707, 80, 18000, 18000, 18000, "`$/caja 15 kilos granel", "Perú", 1200, 15),
    @(150, 44901, 15, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(151, 45016, 50, 25000, 26000, 25600, "`$/malla 20 kilos", "Perú", 1280, 20),
    @(152, 44753, 80, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(153, 44607, 40, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(154, 44680, 20, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(155, 44959, 40, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(156, 44809, 50, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(157, 44452, 50, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(158, 44767, 80, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(159, 44567, 50, 18000, 18000, 18000, "`$/malla 20 kilos", "Región de Arica y Parinacota", 900, 20),
    @(160, 45043, 80, 26000, 26000, 26000, "`$/caja 18 kilos", "Perú", 1444, 18),
    @(161, 44813, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(162, 44496, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(163, 44970, 30, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(164, 44776, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(165, 44389, 45, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(166, 44706, 50, 18000, 18000, 18000, "`$/caja 15 kilos granel", "Perú", 1200, 15),
    @(167, 44972, 40, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(168, 44832, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(169, 44987, 50, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(170, 44762, 15, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(171, 45028, 40, 25000, 25000, 25000, "`$/malla 20 kilos", "Perú", 1250, 20),
    @(172, 44705, 20, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(173, 44690, 15, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(174, 44827, 40, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(175, 45042, 90, 26000, 26000, 26000, "`$/caja 18 kilos", "Perú", 1444, 18),
    @(176, 44497, 30, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(177, 44497, 40, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(178, 44441, 40, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(179, 44587, 55, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(180, 44636, 50, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(181, 44455, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(182, 44669, 25, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(183, 45027, 30, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(184, 45001, 50, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(185, 44795, 50, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(186, 44915, 20, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(187, 44757, 40, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(188, 44817, 25, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(189, 44880, 30, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(190, 44585, 20, 16000, 16000, 16000, "`$/malla 20 kilos", "Perú", 800, 20),
    @(191, 44532, 40, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(192, 44957, 20, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(193, 44340, 40, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(194, 44859, 45, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(195, 44803, 15, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(196, 44341, 40, 17000, 18000, 17500, "`$/malla 20 kilos", "Perú", 875, 20),
    @(197, 44722, 20, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(198, 44741, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(199, 44921, 15, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20),
    @(200, 45014, 40, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(201, 44648, 30, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(202, 44999, 40, 26000, 26000, 26000, "`$/malla 20 kilos", "Perú", 1300, 20),
    @(203, 44662, 20, 18000, 18000, 18000, "`$/caja 15 kilos granel", "Perú", 1200, 15),
    @(204, 44662, 50, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(205, 44910, 50, 24000, 24000, 24000, "`$/malla 20 kilos", "Perú", 1200, 20)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = [double]$row[1]
    $ws.Cells.Item($r, 10).Value = [double]$row[2]
    $ws.Cells.Item($r, 11).Value = [double]$row[3]
    $ws.Cells.Item($r, 12).Value = [double]$row[4]
    $ws.Cells.Item($r, 13).Value = [double]$row[5]
    $ws.Cells.Item($r, 14).Value = $row[6]
    $ws.Cells.Item($r, 15).Value = $row[7]
    $ws.Cells.Item($r, 16).Value = [double]$row[8]
    $ws.Cells.Item($r, 17).Value = [double]$row[9]
}

# Row 205 is brand new: fill in the constant columns that are identical across the whole block
# (A Mercado ID, B Mercado, C Region, E Codreg, F Categoria ID, G Categoria, H Variedad, I Calidad, R Clasificacion)
$ws.Cells.Item(205, 1).Value  = 10
$ws.Cells.Item(205, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(205, 3).Value  = "La Araucanía"
$ws.Cells.Item(205, 5).Value  = 9
$ws.Cells.Item(205, 6).Value  = 100114002
$ws.Cells.Item(205, 7).Value  = "Camote"
$ws.Cells.Item(205, 8).Value  = "Sin especificar"
$ws.Cells.Item(205, 9).Value  = "Primera"
$ws.Cells.Item(205, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D
$ws.Cells.Item(205, 4).NumberFormat = $ws.Cells.Item(204, 4).NumberFormat()
